$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (strikeouts) column values, replacing the old "Strike#" derived values.
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 0
    8  = 2
    9  = 1
    10 = 1
    11 = 3
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 1
    17 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
